$d = $word.ActiveDocument

$d.Content.Find.Execute("RR000000022MA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "RR000000008MA", 2)

$d.Content.Find.Execute("MR.qw qw", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MR.qwe qwe", 2)

$d.Content.Find.Execute("Ain Chegga : qw", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ait Iaaz : qwe", 2)

$d.Content.Find.Execute("12331", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0642187694", 2)
